$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.917.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.007.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.09%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.33%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.007.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.19%  "

$ws.Range("E12").Value = "  -2.33%  "

$ws.Range("E13").Value = "  -4.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.25%  "

$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.499.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.889.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.007.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("E25").Value = "  -6.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.66%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -2.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.29%  "

$ws.Range("E32").Value = "  -6.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("E35").Value = "  -4.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0785"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.06%  "

$ws.Range("E38").Value = "  -5.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.01%  "

$ws.Range("E42").Value = "  -1.09%  "

$ws.Range("E43").Value = "  -6.53%  "

$ws.Range("E44").Value = "  -3.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "377.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -14.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.742.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "

$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("E51").Value = "  -3.41%  "
